$wb = $excel.ActiveWorkbook

# --- "Change History" sheet: log the change ---
$history = $wb.Worksheets.Item("Change History")
$history.Range("B10").Value = (Get-Date -Year 2018 -Month 5 -Day 4).Date
$history.Range("C10").Value = "tfs10890 - ecl disable pilot survey question`nset active to No and end date to 05/04/2018"
$history.Range("D10").Value = 1.05
$history.Range("E10").Value = "Doug Stearns"
$history.Rows.Item(10).RowHeight = 28.8
$history.Range("C10").Select()

# --- "eCL Survey" sheet: disable the pilot survey question on row 7 ---
# (done last / selected last so this sheet - which was already the active
# tab in the original workbook - remains the active tab on save)
$survey = $wb.Worksheets.Item("eCL Survey")
$survey.Range("A7").Value = "No"

# Copy the existing date-formatted style from the "Start Date" column (L)
# onto the new "End Date" entry (M7) before writing its value, so it picks
# up the same plain date number format instead of Excel inventing a new one.
$survey.Range("L2").Copy()
$survey.Range("M7").PasteSpecial(-4122) # xlPasteFormats
$survey.Range("M7").Value = (Get-Date -Year 2018 -Month 5 -Day 4).Date
$survey.Range("A7").Select()
